$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    3.977959584536734,
    7.622334346345951,
    8.073032840776023,
    9.379189874814472,
    2.736651432341773,
    9.799971428978751,
    6.595681556606564,
    5.560753044758911,
    7.876287900559333,
    8.181062382652982,
    1.446861457488097,
    5.792735607280207,
    4.063784253847302,
    2.792715961972817,
    2.431485653114777,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144,
    5.813051703908144
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}
